$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 ("2021年") appended after the existing last data row (11, "2020年").
$ws.Range("A12").Value = "2021年"

# Copy row 11's label-cell formatting (bold/border/center style) onto A12 so it
# matches the style used by every other year label in column A.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B12").Value = 8888
$ws.Range("D12").Value = 6856
$ws.Range("E12").Value = 21525
$ws.Range("F12").Value = 15962
$ws.Range("G12").Value = 44422
$ws.Range("H12").Value = 49652
$ws.Range("I12").Value = 4982
$ws.Range("J12").Value = 2823
$ws.Range("K12").Value = 8261
$ws.Range("M12").Value = 110453
$ws.Range("N12").Value = 27897
$ws.Range("O12").Value = 14038
$ws.Range("P12").Value = 3711
$ws.Range("Q12").Value = 2208
$ws.Range("T12").Value = 1086
$ws.Range("U12").Value = 224577

# C12 / L12 / R12 / S12 stay blank for 2021, same as the blank cells in the
# preceding rows (2017年, 2018年, 2019年, 2020年): they are still present as
# empty-text cells rather than fully absent, so enter an empty string via the
# leading-apostrophe text marker (a plain "" assignment clears the cell
# entirely) and then strip the resulting quote-prefix formatting.
foreach ($col in @("C", "L", "R", "S")) {
    $cell = $ws.Range("$col`12")
    $cell.Value = "'"
    $cell.ClearFormats()
}
